$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("X10").Value = -0.16000300000000323
$ws.Range("Y10").Value = "Down"

$ws.Range("A11").NumberFormat = "m/d/yy h:mm"
$ws.Range("S11").NumberFormat = $ws.Range("S10").NumberFormat
$ws.Range("T11").NumberFormat = $ws.Range("T10").NumberFormat

$ws.Range("A11").Value = 42654.894490740742
$ws.Range("B11").Value = 14
$ws.Range("C11").Value = "Buy"
$ws.Range("D11").Value = 10
$ws.Range("E11").Value = 15067
$ws.Range("F11").Value = 891
$ws.Range("G11").Value = 52
$ws.Range("H11").Value = 46
$ws.Range("I11").Value = 79
$ws.Range("J11").Value = 19
$ws.Range("K11").Value = 21085
$ws.Range("L11").Value = 143
$ws.Range("M11").Value = 127
$ws.Range("N11").Value = 44
$ws.Range("O11").Value = 11
$ws.Range("P11").Value = "Named"
$ws.Range("Q11").Value = 28.689659976213832
$ws.Range("R11").Value = 0.84
$ws.Range("S11").Value = -0.0125
$ws.Range("T11").Value = -0.0261
$ws.Range("U11").Value = 14.56
$ws.Range("V11").Value = "N/A"
$ws.Range("W11").Value = 0
